$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update fuel type names to uppercase
$ws.Range("A2").Value = "PETROL"
$ws.Range("A3").Value = "DIESEL"

# Update selection to A3 as last action
$ws.Range("A3").Select()
